$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "elizabeth" user record in row 9
$ws.Range("A9").Value = "elizabeth"

# B9 gets the email text plus a mailto hyperlink, matching the style of B2:B8
$ws.Range("B9").Value = "liza@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:liza@gmail.com")
$ws.Range("B9").Style = $ws.Range("B8").Style

$ws.Range("C9").Value = "liza1234"
$ws.Range("D9").Value = "b.png"

# Leave the selection on the newly added row, as in the source edit
$null = $ws.Range("D8:D9").Select()

$wb.Save()
